# Lab05 Group A instructions: clarify the upload-count wording, split the
# "Two files" / "Four files" bullet descriptions into itemized run pieces,
# merge the code-review sentence back together, and relocate the Word
# "_GoBack" last-edit bookmark so it spans the whole block that was touched
# (start of the "Upload the following..." paragraph through the end of the
# "...filled in by you." paragraph) -- matching how Word itself stamps
# _GoBack after a multi-paragraph editing session.

$d = $word.ActiveDocument

$searchRange = $d.Content
$found = $searchRange.Find.Execute("Upload the following to the", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'Upload the following to the' paragraph"
}

$p1 = $searchRange.Paragraphs(1)
$p2 = $p1.Next()
$p3 = $p2.Next()
$p4 = $p3.Next()

$block = $d.Range($p1.Range.Start, $p4.Range.End)

$blockXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="301C6E56" w14:textId="086153D9" w:rsidR="001D214A" w:rsidRPr="00351FA2" w:rsidRDefault="001D214A" w:rsidP="007F033B"><w:pPr><w:spacing w:before="120" w:after="120"/><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Upload the following</w:t></w:r><w:r w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> 7 files</w:t></w:r><w:r w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> to the </w:t></w:r><w:r w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Lab Production Version</w:t></w:r><w:r w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> assignment:</w:t></w:r></w:p><w:p w14:paraId="4DC19B07" w14:textId="16738BDB" w:rsidR="001D214A" w:rsidRPr="00351FA2" w:rsidRDefault="001D214A" w:rsidP="007F033B"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:spacing w:before="120" w:after="120"/><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>T</w:t></w:r><w:r w:rsidR="00B2055F"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>wo files (.html and .</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00B2055F"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>js</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00B2055F"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>) for</w:t></w:r><w:r w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> part 1.</w:t></w:r></w:p><w:p w14:paraId="176547D1" w14:textId="3A9B7568" w:rsidR="001D214A" w:rsidRPr="00351FA2" w:rsidRDefault="00B2055F" w:rsidP="007F033B"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:spacing w:before="120" w:after="120"/><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Four</w:t></w:r><w:r w:rsidR="001D214A" w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> files</w:t></w:r><w:r w:rsidR="001D214A" w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="001D214A" w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>(</w:t></w:r><w:r w:rsidR="001D214A" w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">2 </w:t></w:r><w:r w:rsidR="001D214A" w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">html and </w:t></w:r><w:r w:rsidR="001D214A" w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">2 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D214A" w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>j</w:t></w:r><w:r w:rsidR="001D214A" w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D214A" w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">) </w:t></w:r><w:r w:rsidR="001D214A" w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>for part 2.</w:t></w:r></w:p><w:p w14:paraId="01E598C1" w14:textId="1BFBF209" w:rsidR="005F66B6" w:rsidRPr="00351FA2" w:rsidRDefault="001D214A" w:rsidP="007F033B"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:spacing w:before="120" w:after="120"/><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">The code review </w:t></w:r><w:r w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/></w:rPr><w:t>from</w:t></w:r><w:r w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> your lab partner with the </w:t></w:r><w:r w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/></w:rPr><w:t>“Prod” column filled in by you</w:t></w:r><w:r w:rsidRPr="00351FA2"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>.</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$block.InsertXML($blockXml)
